# "update all templates and sample files in inst #48"
#
# On the "Instructions" sheet, two red "last updated" notices are added
# in column C next to the first two instruction bullets:
#   C1 -> "Template updated 12/8/22."
#   C2 -> "Samples updated 1/8/23"
# Both use the default (Calibri 11) font, colored red, with no other
# special formatting (no border/fill/alignment changes).

$wb = $excel.ActiveWorkbook

$wsInstructions = $wb.Worksheets.Item("Instructions")

$wsInstructions.Range("C1").Value = "Template updated 12/8/22."
$wsInstructions.Range("C1").Font.Color = 255

$wsInstructions.Range("C2").Value = "Samples updated 1/8/23"
$wsInstructions.Range("C2").Font.Color = 255

# Reflect the cursor/selection state left behind in the saved file: the
# last active cell on the Instructions sheet is C3 (right below the new
# notes) and on the Completeness sheet is D12.
[void]$wsInstructions.Range("C3").Select()

$wsCompleteness = $wb.Worksheets.Item("Completeness")
[void]$wsCompleteness.Range("D12").Select()
